# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the "Maduin_Profits" workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# ALC row 43
$ws.Range("H43").Value = 4960.231
$ws.Range("I43").Value = 3666.6667
$ws.Range("J43").Value = 6069
$ws.Range("K43").Value = 3666.6667
$ws.Range("L43").Value = 6069
$ws.Range("M43").Value = -3597.6667
$ws.Range("N43").Value = -6207

# ALC row 62
$ws.Range("H62").Value = 3462.6667
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126

# ALC row 65
$ws.Range("H65").Value = 3462.6667
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630

# ALC row 70
$ws.Range("H70").Value = 2355.375
$ws.Range("I70").Value = 2140.5
$ws.Range("K70").Value = 6421.5
$ws.Range("M70").Value = -6151.5

# ALC row 73
$ws.Range("H73").Value = 2355.375
$ws.Range("I73").Value = 2140.5
$ws.Range("K73").Value = 6421.5
$ws.Range("M73").Value = -5485.5

$ws = $wb.Worksheets.Item("ARM")

# ARM row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# ARM row 45
$ws.Range("H45").Value = 4018.3076
$ws.Range("I45").Value = 1450
$ws.Range("J45").Value = 4788.8
$ws.Range("K45").Value = 1450
$ws.Range("L45").Value = 4788.8
$ws.Range("M45").Value = -1073
$ws.Range("N45").Value = -5542.8

# ARM row 47
$ws.Range("H47").Value = 25000
$ws.Range("J47").Value = 25000
$ws.Range("L47").Value = 25000
$ws.Range("N47").Value = -26450

# ARM row 61
$ws.Range("H61").Value = 1360.125
$ws.Range("I61").Value = 1360.125
$ws.Range("K61").Value = 1360.125
$ws.Range("M61").Value = -1148.125

# ARM row 74
$ws.Range("H74").Value = 948
$ws.Range("I74").Value = 948
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 948
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -74
$ws.Range("N74").ClearContents()

# ARM row 77
$ws.Range("H77").Value = 948
$ws.Range("I77").Value = 948
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4740
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -372
$ws.Range("N77").ClearContents()

# ARM row 88
$ws.Range("H88").Value = 3350
$ws.Range("I88").Value = 3585
$ws.Range("J88").Value = 3249.2856
$ws.Range("K88").Value = 3585
$ws.Range("L88").Value = 3249.2856
$ws.Range("M88").Value = -3179
$ws.Range("N88").Value = -4061.2856

# ARM row 91
$ws.Range("H91").Value = 3350
$ws.Range("I91").Value = 3585
$ws.Range("J91").Value = 3249.2856
$ws.Range("K91").Value = 3585
$ws.Range("L91").Value = 3249.2856
$ws.Range("M91").Value = -2181
$ws.Range("N91").Value = -6057.2856

# ARM row 97
$ws.Range("H97").Value = 1974.1111
$ws.Range("I97").Value = 801.2
$ws.Range("K97").Value = 801.2
$ws.Range("M97").Value = -305.2

# ARM row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ARM row 132
$ws.Range("H132").Value = 885.61536
$ws.Range("I132").Value = 885.61536
$ws.Range("K132").Value = 2656.84608
$ws.Range("M132").Value = -126.8460800000003

# ARM row 136
$ws.Range("H136").Value = 1360.125
$ws.Range("I136").Value = 1360.125
$ws.Range("K136").Value = 4080.375
$ws.Range("M136").Value = -1530.375

$ws = $wb.Worksheets.Item("BSM")

# BSM row 20
$ws.Range("H20").Value = 1549
$ws.Range("I20").Value = 828.7143
$ws.Range("K20").Value = 828.7143
$ws.Range("M20").Value = -581.7143

# BSM row 86
$ws.Range("H86").Value = 5017
$ws.Range("I86").Value = 4651.1665
$ws.Range("J86").Value = 5382.8335
$ws.Range("K86").Value = 4651.1665
$ws.Range("L86").Value = 5382.8335
$ws.Range("M86").Value = -3528.1665
$ws.Range("N86").Value = -7628.8335

# BSM row 89
$ws.Range("H89").Value = 5017
$ws.Range("I89").Value = 4651.1665
$ws.Range("J89").Value = 5382.8335
$ws.Range("K89").Value = 23255.8325
$ws.Range("L89").Value = 26914.1675
$ws.Range("M89").Value = -17639.8325
$ws.Range("N89").Value = -38146.1675

# BSM row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CRP")

# CRP row 31
$ws.Range("H31").Value = 2605
$ws.Range("I31").Value = 1045.6666
$ws.Range("J31").Value = 4164.3335
$ws.Range("K31").Value = 1045.6666
$ws.Range("L31").Value = 4164.3335
$ws.Range("M31").Value = -750.6666
$ws.Range("N31").Value = -4754.3335

# CRP row 34
$ws.Range("H34").Value = 2605
$ws.Range("I34").Value = 1045.6666
$ws.Range("J34").Value = 4164.3335
$ws.Range("K34").Value = 1045.6666
$ws.Range("L34").Value = 4164.3335
$ws.Range("M34").Value = -843.6666
$ws.Range("N34").Value = -4568.3335

# CRP row 35
$ws.Range("H35").Value = 1653.5
$ws.Range("I35").Value = 1653.5
$ws.Range("K35").Value = 1653.5
$ws.Range("M35").Value = -1359.5

# CRP row 122
$ws.Range("H122").Value = 1979.75
$ws.Range("I122").Value = 1796.091
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5388.272999999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2938.272999999999
$ws.Range("N122").Value = -16900

# CRP row 134
$ws.Range("H134").Value = 1170
$ws.Range("I134").Value = 1170
$ws.Range("K134").Value = 3510
$ws.Range("M134").Value = -975

$ws = $wb.Worksheets.Item("CUL")

# CUL row 10
$ws.Range("H10").Value = 27.666666
$ws.Range("I10").Value = 27.666666
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 82.99999800000001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 56.00000199999999
$ws.Range("N10").ClearContents()

# CUL row 22
$ws.Range("H22").Value = 14710.1
$ws.Range("I22").Value = 17613.334
$ws.Range("K22").Value = 52840.00199999999
$ws.Range("M22").Value = -52671.00199999999

# CUL row 27
$ws.Range("H27").Value = 14710.1
$ws.Range("I27").Value = 17613.334
$ws.Range("K27").Value = 52840.00199999999
$ws.Range("M27").Value = -52738.00199999999

# CUL row 41
$ws.Range("H41").Value = 139.33333
$ws.Range("I41").Value = 139.33333
$ws.Range("K41").Value = 417.99999
$ws.Range("M41").Value = -79.99998999999997

# CUL row 47
$ws.Range("H47").Value = 169.5
$ws.Range("I47").Value = 192.66667
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 578.00001
$ws.Range("L47").Value = 300
$ws.Range("M47").Value = -147.00001
$ws.Range("N47").Value = -1162

# CUL row 75
$ws.Range("H75").Value = 1015
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1015
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 3045
$ws.Range("N75").Value = -5041
$ws.Range("M75").ClearContents()

# CUL row 78
$ws.Range("H78").Value = 1015
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1015
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 9135
$ws.Range("N78").Value = -19119
$ws.Range("M78").ClearContents()

# CUL row 93
$ws.Range("H93").Value = 2730.1428
$ws.Range("I93").Value = 222.6
$ws.Range("J93").Value = 8999
$ws.Range("K93").Value = 667.8
$ws.Range("L93").Value = 26997
$ws.Range("M93").Value = 1204.2
$ws.Range("N93").Value = -30741

# CUL row 138
$ws.Range("H138").Value = 2161.5
$ws.Range("I138").Value = 1742.5
$ws.Range("J138").Value = 2999.5
$ws.Range("K138").Value = 5227.5
$ws.Range("L138").Value = 8998.5
$ws.Range("M138").Value = -87.5
$ws.Range("N138").Value = -19278.5

$ws = $wb.Worksheets.Item("GSM")

# GSM row 70
$ws.Range("H70").Value = 12000.25
$ws.Range("I70").Value = 12000.25
$ws.Range("K70").Value = 12000.25
$ws.Range("M70").Value = -11730.25

# GSM row 73
$ws.Range("H73").Value = 12000.25
$ws.Range("I73").Value = 12000.25
$ws.Range("K73").Value = 12000.25
$ws.Range("M73").Value = -11064.25

# GSM row 132
$ws.Range("H132").Value = 2315.6155
$ws.Range("I132").Value = 2315.6155
$ws.Range("K132").Value = 6946.8465
$ws.Range("M132").Value = -4416.8465

$ws = $wb.Worksheets.Item("LTW")

# LTW row 9
$ws.Range("H9").Value = 391.5
$ws.Range("I9").Value = 447.6
$ws.Range("K9").Value = 447.6
$ws.Range("M9").Value = -223.6

# LTW row 46
$ws.Range("H46").Value = 5090.909
$ws.Range("I46").Value = 4600
$ws.Range("K46").Value = 4600
$ws.Range("M46").Value = -4412

# LTW row 68
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2251
$ws.Range("N68").ClearContents()

# LTW row 71
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -11256
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")

# WVR row 96
$ws.Range("H96").Value = 1016.625
$ws.Range("I96").Value = 823.25
$ws.Range("J96").Value = 1210
$ws.Range("K96").Value = 823.25
$ws.Range("L96").Value = 1210
$ws.Range("M96").Value = 549.75
$ws.Range("N96").Value = -3956

